$d = $word.ActiveDocument

$replacements = @(
    @("15×61=", "15×64="),
    @("53×91=", "84×34="),
    @("78×44=", "19×76="),
    @("69×58=", "60×91="),
    @("11×67=", "47×21="),
    @("32×21=", "70×51="),
    @("25×68=", "21×76="),
    @("14×57=", "66×80="),
    @("11×42=", "56×91="),
    @("37×53=", "25×39="),
    @("51×92=", "81×56="),
    @("37×84=", "52×78="),
    @("96×92=", "87×47="),
    @("15×68=", "36×63="),
    @("33×50=", "50×59="),
    @("88×52=", "40×21="),
    @("85×96=", "22×52="),
    @("64×13=", "23×90="),
    @("33×17=", "91×16="),
    @("46×49=", "16×56="),
    @("80×65=", "34×89="),
    @("59×25=", "60×86="),
    @("68×60=", "78×77="),
    @("54×47=", "71×54="),
    @("37×74=", "64×36=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
